$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data to append (row 234 through 238), matching the
# existing pattern: column A = date serial (styled like existing col A),
# B = nuovi pos., C = somma mobile 7gg., D = somma mobile 7gg. per 100mila abitanti

$newRows = @(
    @{ Row = 234; A = 44308; B = 3; C = 5; D = 109.051254089422 },
    @{ Row = 235; A = 44309; B = 4; C = 8; D = 174.4820065430752 },
    @{ Row = 236; A = 44310; B = 0; C = 8; D = 174.4820065430752 },
    @{ Row = 237; A = 44311; B = 0; C = 7; D = 152.6717557251908 },
    @{ Row = 238; A = 44312; B = 0; C = 7; D = 152.6717557251908 }
)

foreach ($r in $newRows) {
    $rowIndex = $r.Row

    # Column A: date value, copy style from the cell above so it keeps
    # the same date formatting/style as the rest of the column.
    $ws.Cells.Item($rowIndex, 1).Value = $r.A
    $ws.Cells.Item($rowIndex - 1, 1).Copy()
    $ws.Cells.Item($rowIndex, 1).PasteSpecial(-4122)

    $ws.Cells.Item($rowIndex, 1).Value = $r.A
    $ws.Cells.Item($rowIndex, 2).Value = $r.B
    $ws.Cells.Item($rowIndex, 3).Value = $r.C
    $ws.Cells.Item($rowIndex, 4).Value = $r.D
}
